$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Remove the "Power BI Administrator" role row entirely (row 21).
#    Everything below shifts up by one row.
$ws.Rows.Item(21).Delete()

# 2) "Helpdesk Administrator" (now at row 56, after the shift above) moves
#    into the alphabetically-sorted "eligible" block, landing right after
#    "User Administrator" (row 29) at row 30, and gains an Eligible value.
#    Insert a blank row at 30 (pushes "Application Developer" and below
#    down by one) and fill it in, then delete the now-duplicate old row.
$ws.Rows.Item(30).Insert()
$ws.Range("A30").Value = "Helpdesk Administrator"
$ws.Range("C30").Value = "any.admin@alyaconsulting.ch"
$ws.Rows.Item(57).Delete()

# 3) Append the three newly-introduced roles at the bottom of the list.
$ws.Range("A101").Value = "Fabric Administrator"
$ws.Range("A102").Value = "Global Secure Access Administrator"
$ws.Range("A103").Value = "Extended Directory User Administrator"
